# Apply the target edit:
#  1. Rename the sheet from "UniformA-HW30.xpc" to "UniformA"
#  2. Append a new row 16 to the sheet that mirrors row 15's layout
#     (index 14, same HexGrid-60degTilt5degRes label, and 1's across C:P)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet/tab.
$ws.Name = "UniformA"

# 2. Populate the new row of data.
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "HexGrid-60degTilt5degRes"
for ($col = 3; $col -le 16; $col++) {
    $ws.Cells.Item(16, $col).Value = 1
}

# Mirror the formatting of column A's label cell (row 15) onto the new
# row 15 -> 16 entry so the same cell style (bold / bordered / centered)
# is reused instead of a brand new style being created.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

Write-Output "Renamed sheet and appended row 16"
